{"js": "// The document ends with two empty paragraphs right before the final\n// section break (after the last table). The second (last) of those two\n// empty paragraphs becomes \"Rare Patterns:\" and a brand new paragraph\n// with the pattern-count data is added right after it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"Rare Patterns:\", \"Replace\");\nlastParagraph.insertParagraph(\n  \"{1:2, 2:3, 3:3, 4:1, 5:3, {1,3}:2, {1,4}:1, {1,2}:1, {1,5}:1, \",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# The document ends with two empty paragraphs right before the final\n# section break (after the last table). The second (last) of those two\n# empty paragraphs becomes \"Rare Patterns:\" and a brand new paragraph\n# with the pattern-count data is added right after it.\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertBefore(\"Rare Patterns:\")\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.InsertBefore(\"{1:2, 2:3, 3:3, 4:1, 5:3, {1,3}:2, {1,4}:1, {1,2}:1, {1,5}:1, \")\n"}
